$d = $word.ActiveDocument

# --- Change 1: merge the runs of the "Mein Leben dreht sich um Lernen..." paragraph
# into a single run with the identical text (removes the stray w:proofErr markers and
# the mid-sentence run splits left over from spell/grammar-check squiggles).
$mergedText = "Mein Leben dreht sich um Lernen, Arbeiten(man muss Geld verdienen) und Schlafen. Eine franzözische Redewendung dafür ist:Metro-Boulot-Dodo."
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# --- Change 2: shorten the closing "Mit freundlichen Grüßen" to "MFG"
$d.Content.Find.Execute("Mit freundlichen Grüßen", $true, $false, $false, $false, $false, $true, 1, $false, "MFG", 2) | Out-Null
